$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated financial figures per cell (row 2-9, columns D..AJ)
$values = @{
    "D2" = 4581
    "E2" = 270
    "F2" = 270
    "G2" = 207
    "H2" = 153
    "I2" = 153
    "K2" = 5442
    "L2" = 3030
    "M2" = 2413
    "N2" = 2413
    "P2" = 102
    "Q2" = 267
    "R2" = -156
    "S2" = -97
    "T2" = 133
    "U2" = 134
    "V2" = 1477
    "W2" = 5.89
    "X2" = 3.35
    "Y2" = 6.5
    "Z2" = 2.86
    "AA2" = 125.58
    "AB2" = 2263.83
    "AC2" = 750
    "AD2" = 12.07
    "AE2" = 11790
    "AF2" = 0.77
    "AG2" = 200
    "AH2" = 2.21
    "AI2" = 26.68
    "AJ2" = 20000000
    "D3" = 4616
    "E3" = 307
    "F3" = 307
    "G3" = 421
    "H3" = 327
    "I3" = 327
    "K3" = 5729
    "L3" = 3017
    "M3" = 2712
    "N3" = 2712
    "P3" = 102
    "Q3" = 289
    "R3" = -92
    "S3" = -138
    "T3" = 98
    "U3" = 191
    "V3" = 1379
    "W3" = 6.64
    "X3" = 7.09
    "Y3" = 12.78
    "Z3" = 5.86
    "AA3" = 111.24
    "AB3" = 2548.03
    "AC3" = 1600
    "AD3" = 6.97
    "AE3" = 13253
    "AF3" = 0.84
    "AG3" = 200
    "AH3" = 1.79
    "AI3" = 12.51
    "AJ3" = 20000000
    "D4" = 4789
    "E4" = 312
    "F4" = 312
    "G4" = 207
    "H4" = 95
    "I4" = 95
    "K4" = 5706
    "L4" = 2941
    "M4" = 2765
    "N4" = 2765
    "P4" = 102
    "Q4" = 445
    "R4" = -258
    "S4" = -139
    "T4" = 102
    "U4" = 343
    "V4" = 1282
    "W4" = 6.51
    "X4" = 1.98
    "Y4" = 3.47
    "Z4" = 1.66
    "AA4" = 106.39
    "AB4" = 2597.66
    "AC4" = 464
    "AD4" = 18.42
    "AE4" = 13510
    "AF4" = 0.63
    "AG4" = 200
    "AH4" = 2.34
    "AI4" = 43.1
    "AJ4" = 20000000
    "D5" = 5514
    "E5" = 278
    "F5" = 278
    "G5" = 398
    "H5" = 294
    "I5" = 294
    "K5" = 5918
    "L5" = 2928
    "M5" = 2990
    "N5" = 2990
    "P5" = 102
    "Q5" = 333
    "R5" = -172
    "S5" = -263
    "T5" = 121
    "U5" = 212
    "V5" = 1192
    "W5" = 5.03
    "X5" = 5.34
    "Y5" = 10.23
    "Z5" = 5.07
    "AA5" = 97.91
    "AB5" = 2844.76
    "AC5" = 1439
    "AD5" = 6.03
    "AE5" = 14613
    "AF5" = 0.59
    "AG5" = 250
    "AH5" = 2.88
    "AI5" = 17.38
    "AJ5" = 20000000
    "D6" = 6147
    "E6" = 227
    "F6" = 227
    "G6" = 191
    "H6" = 127
    "I6" = 127
    "K6" = 6074
    "L6" = 2996
    "M6" = 3078
    "N6" = 3078
    "P6" = 102
    "Q6" = 140
    "R6" = -124
    "S6" = -29
    "T6" = 98
    "U6" = 42
    "V6" = 1216
    "W6" = 3.69
    "X6" = 2.07
    "Y6" = 4.2
    "Z6" = 2.13
    "AA6" = 97.34999999999999
    "AB6" = 2917.21
    "AC6" = 623
    "AD6" = 13.22
    "AE6" = 15042
    "AF6" = 0.55
    "AG6" = 250
    "AH6" = 3.04
    "AI6" = 40.17
    "AJ6" = 20000000
    "D7" = 6559
    "E7" = 376
    "G7" = 373
    "H7" = 265
    "I7" = 265
    "K7" = 6283
    "L7" = 2987
    "M7" = 3296
    "N7" = 3296
    "P7" = 101
    "Q7" = 470
    "R7" = -87
    "S7" = -314
    "T7" = 75
    "W7" = 5.73
    "X7" = 4.04
    "Y7" = 8.31
    "Z7" = 4.29
    "AA7" = 90.63
    "AC7" = 1295
    "AD7" = 6.06
    "AE7" = 16107
    "AF7" = 0.49
    "AG7" = 300
    "AH7" = 3.82
    "AI7" = 22.64
    "D8" = 6888
    "E8" = 384
    "G8" = 374
    "H8" = 278
    "I8" = 278
    "K8" = 6548
    "L8" = 3022
    "M8" = 3426
    "N8" = 3526
    "P8" = 101
    "Q8" = 346
    "R8" = -98
    "S8" = -92
    "T8" = 76
    "W8" = 5.58
    "X8" = 4.04
    "Y8" = 8.15
    "Z8" = 4.33
    "AA8" = 88.20999999999999
    "AC8" = 1359
    "AD8" = 5.45
    "AE8" = 17229
    "AF8" = 0.43
    "AG8" = 300
    "AH8" = 4.05
    "AI8" = 21.58
    "D9" = 7096
    "E9" = 422
    "G9" = 397
    "H9" = 300
    "I9" = 300
    "K9" = 6805
    "L9" = 3032
    "M9" = 3768
    "N9" = 3768
    "P9" = 101
    "Q9" = 373
    "R9" = -188
    "S9" = -100
    "T9" = 171
    "W9" = 5.94
    "X9" = 4.23
    "Y9" = 8.23
    "Z9" = 4.49
    "AA9" = 80.48999999999999
    "AC9" = 1466
    "AD9" = 5.05
    "AE9" = 18411
    "AF9" = 0.4
    "AG9" = 350
    "AH9" = 23.33
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# Cells removed entirely in the updated report (columns no longer reported)
$clearedCells = @(
    "J2",
    "O2",
    "J3",
    "O3",
    "J4",
    "O4",
    "J5",
    "O5",
    "U7",
    "U8",
    "U9",
)

foreach ($ref in $clearedCells) {
    $ws.Range($ref).ClearContents()
}

Write-Host "Applied IFRS list corrections"